$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "64.378.58"
$ws.Range("E2").Value2 = "  -2.22%  "

$ws.Range("D3").Value2 = "3.205.69"
$ws.Range("E3").Value2 = "  -7.23%  "

$ws.Range("D5").Value = "'560.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -3.98%  "

$ws.Range("D6").Value = "'173.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +0.84%  "

$ws.Range("B7").Value2 = "USDC"
$ws.Range("C7").Value2 = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  -0.04%  "

$ws.Range("B8").Value2 = "XRP"
$ws.Range("C8").Value2 = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +0.49%  "

$ws.Range("D9").Value2 = "3.207.15"
$ws.Range("E9").Value2 = "  -7.10%  "

$ws.Range("D10").Value = "'0.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -5.43%  "

$ws.Range("D11").Value = "'6.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -2.49%  "

$ws.Range("D12").Value = "'0.395"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -3.23%  "

$ws.Range("D13").Value2 = "3.769.69"
$ws.Range("E13").Value2 = "  -7.22%  "

$ws.Range("D14").Value = "'0.134"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -0.22%  "

$ws.Range("D15").Value = "'27.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -5.23%  "

$ws.Range("D16").Value2 = "64.577.16"
$ws.Range("E16").Value2 = "  -2.06%  "

$ws.Range("D17").Value = "'0.0000162"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -4.76%  "

$ws.Range("D18").Value2 = "3.205.43"
$ws.Range("E18").Value2 = "  -7.58%  "

$ws.Range("D19").Value = "'5.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -4.83%  "

$ws.Range("D20").Value = "'12.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -6.32%  "

$ws.Range("D21").Value = "'355.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -3.04%  "

$ws.Range("D22").Value = "'7.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -6.11%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -0.11%  "

$ws.Range("D24").Value = "'69.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -4.98%  "

$ws.Range("D25").Value = "'0.0000119"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -3.49%  "

$ws.Range("D26").Value = "'0.501"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -5.85%  "

$ws.Range("D27").Value = "'9.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -1.78%  "

$ws.Range("D28").Value = "'0.175"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -1.61%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -0.14%  "

$ws.Range("B30").Value2 = "USDe"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +0.04%  "

$ws.Range("D31").Value = "'5.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -2.18%  "

$ws.Range("B32").Value2 = "PancakeSwap"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -4.52%  "

$ws.Range("D33").Value = "'22.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -7.08%  "

$ws.Range("D34").Value = "'6.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -6.51%  "

$ws.Range("D35").Value = "'1.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -7.70%  "

$ws.Range("B36").Value2 = "Monero"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'158.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -0.69%  "

$ws.Range("B37").Value2 = "ImmutableX"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -4.52%  "

$ws.Range("D38").Value = "'0.815"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -8.05%  "

$ws.Range("D39").Value = "'26.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -9.44%  "

$ws.Range("B40").Value2 = "dogwifhat"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'2.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -4.28%  "

$ws.Range("B41").Value2 = "Stacks"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -3.16%  "

$ws.Range("B42").Value2 = "Maker"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value2 = "2.663.23"
$ws.Range("E42").Value2 = "  -4.17%  "

$ws.Range("D43").Value = "'5.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -6.38%  "

$ws.Range("D44").Value = "'4.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -6.89%  "

$ws.Range("D45").Value = "'39.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -1.95%  "

$ws.Range("D46").Value = "'0.0650"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -4.40%  "

$ws.Range("D47").Value = "'321.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -1.85%  "

$ws.Range("D48").Value = "'23.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -2.47%  "

$ws.Range("D49").Value = "'0.0270"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -6.06%  "

$ws.Range("D50").Value = "'0.101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -0.13%  "

$ws.Range("B51").Value2 = "FirstDigitalUSD"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +0.02%  "
